$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").WrapText = $true
$ws.Range("A1").Borders.Item(10).LineStyle = 1
$ws.Range("A1").Borders.Item(10).Weight = 2
$ws.Range("A1").Borders.Item(8).LineStyle = 1
$ws.Range("A1").Borders.Item(8).Weight = 2
$ws.Range("A1").Borders.Item(7).LineStyle = 1
$ws.Range("A1").Borders.Item(7).Weight = 2

$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B1").WrapText = $true
$ws.Range("B1").Borders.Item(10).LineStyle = 1
$ws.Range("B1").Borders.Item(10).Weight = 2
$ws.Range("B1").Borders.Item(8).LineStyle = 1
$ws.Range("B1").Borders.Item(8).Weight = 2
$ws.Range("B1").Borders.Item(7).LineStyle = 1
$ws.Range("B1").Borders.Item(7).Weight = 2
